{"js": "// Renumber the four \"Table N: ...\" headings that precede each results\n// table. The tables/content themselves are unchanged; only the heading\n// numbers move to reflect the new table order.\nconst renumbers = [\n  {\n    oldText: \"Table 4: Maternal Micronutrients and Child Growth Status\",\n    newText: \"Table 5: Maternal Micronutrients and Child Growth Status\",\n  },\n  {\n    oldText: \"Table 2: Maternal Plasma Cortisol and Child Growth Status\",\n    newText: \"Table 3: Maternal Plasma Cortisol and Child Growth Status\",\n  },\n  {\n    oldText: \"Table 5: Maternal Estriol and Child Growth Status\",\n    newText: \"Table 2: Maternal Estriol and Child Growth Status\",\n  },\n  {\n    oldText: \"Table 3: Maternal Inflammation and Child Growth Status\",\n    newText: \"Table 4: Maternal Inflammation and Child Growth Status\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of renumbers) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Renumber the four \"Table N: ...\" headings that precede each results\n# table (Titre2-styled paragraphs). The underlying tables/content are\n# unchanged; only the heading numbers move per the new table order.\n$d = $word.ActiveDocument\n\n$renumbers = @(\n    @{ Old = \"Table 4: Maternal Micronutrients and Child Growth Status\"; New = \"Table 5: Maternal Micronutrients and Child Growth Status\" },\n    @{ Old = \"Table 2: Maternal Plasma Cortisol and Child Growth Status\"; New = \"Table 3: Maternal Plasma Cortisol and Child Growth Status\" },\n    @{ Old = \"Table 5: Maternal Estriol and Child Growth Status\"; New = \"Table 2: Maternal Estriol and Child Growth Status\" },\n    @{ Old = \"Table 3: Maternal Inflammation and Child Growth Status\"; New = \"Table 4: Maternal Inflammation and Child Growth Status\" }\n)\n\nforeach ($item in $renumbers) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $item.Old\n    $find.Replacement.Text = $item.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$item.Old, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$item.New, 2) | Out-Null\n}\n"}
